# Update leve-profit calculation columns (H,I,J,K,L,M,N) across multiple
# job sheets in the workbook, reflecting refreshed market-price data from
# the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 214.8
$ws.Range("I2").Value = 214.8
$ws.Range("K2").Value = 214.8
$ws.Range("M2").Value = -101.8
$ws.Range("H29").Value = 1294.4445
$ws.Range("I29").Value = 1883.3334
$ws.Range("K29").Value = 5650.0002
$ws.Range("M29").Value = -5369.0002
$ws.Range("H38").Value = 6762.125
$ws.Range("J38").Value = 10899.333
$ws.Range("L38").Value = 32697.999
$ws.Range("N38").Value = -33441.999
$ws.Range("H58").Value = 1223
$ws.Range("I58").Value = 91.25
$ws.Range("J58").Value = 5750
$ws.Range("K58").Value = 273.75
$ws.Range("L58").Value = 17250
$ws.Range("M58").Value = -123.75
$ws.Range("N58").Value = -17550
$ws.Range("H76").Value = 4999.3335
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4999.3335
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 4999.3335
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -5629.3335
$ws.Range("H79").Value = 4999.3335
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4999.3335
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 4999.3335
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -7183.3335
$ws.Range("H98").Value = 1549.1154
$ws.Range("I98").Value = 1255.7142
$ws.Range("J98").Value = 2781.4
$ws.Range("K98").Value = 1255.7142
$ws.Range("L98").Value = 2781.4
$ws.Range("M98").Value = 242.2858000000001
$ws.Range("N98").Value = -5777.4
$ws.Range("H122").Value = 1549.1154
$ws.Range("I122").Value = 1255.7142
$ws.Range("J122").Value = 2781.4
$ws.Range("K122").Value = 3767.1426
$ws.Range("L122").Value = 8344.200000000001
$ws.Range("M122").Value = -1317.1426
$ws.Range("N122").Value = -13244.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23981.941
$ws.Range("I32").Value = 23860.182
$ws.Range("K32").Value = 23860.182
$ws.Range("M32").Value = -23573.182
$ws.Range("H61").Value = 3334.0908
$ws.Range("I61").Value = 2460.8708
$ws.Range("K61").Value = 2460.8708
$ws.Range("M61").Value = -2248.8708
$ws.Range("H74").Value = 31251050
$ws.Range("I74").Value = 41667532
$ws.Range("J74").Value = 1600
$ws.Range("K74").Value = 41667532
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = -41666658
$ws.Range("N74").Value = -3348
$ws.Range("H77").Value = 31251050
$ws.Range("I77").Value = 41667532
$ws.Range("J77").Value = 1600
$ws.Range("K77").Value = 208337660
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = -208333292
$ws.Range("N77").Value = -16736
$ws.Range("H102").Value = 980465.1
$ws.Range("I102").Value = 1143624
$ws.Range("J102").Value = 1512
$ws.Range("K102").Value = 1143624
$ws.Range("L102").Value = 1512
$ws.Range("M102").Value = -1142002
$ws.Range("N102").Value = -4756
$ws.Range("H110").Value = 817441.1
$ws.Range("I110").Value = 928633.3
$ws.Range("J110").Value = 2031.6666
$ws.Range("K110").Value = 928633.3
$ws.Range("L110").Value = 2031.6666
$ws.Range("M110").Value = -926588.3
$ws.Range("N110").Value = -6121.6666
$ws.Range("H132").Value = 6291.484
$ws.Range("I132").Value = 2158.8572
$ws.Range("K132").Value = 6476.571599999999
$ws.Range("M132").Value = -3946.571599999999
$ws.Range("H136").Value = 3334.0908
$ws.Range("I136").Value = 2460.8708
$ws.Range("K136").Value = 7382.6124
$ws.Range("M136").Value = -4832.6124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 92779
$ws.Range("J132").Value = 92779
$ws.Range("L132").Value = 92779
$ws.Range("N132").Value = -102899
$ws.Range("H134").Value = 2424.276
$ws.Range("I134").Value = 1731.16
$ws.Range("K134").Value = 5193.48
$ws.Range("M134").Value = -2658.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 919.6111
$ws.Range("I16").Value = 686.5
$ws.Range("J16").Value = 1385.8334
$ws.Range("K16").Value = 686.5
$ws.Range("L16").Value = 1385.8334
$ws.Range("M16").Value = -399.5
$ws.Range("N16").Value = -1959.8334
$ws.Range("H31").Value = 13700560
$ws.Range("I31").Value = 15626338
$ws.Range("K31").Value = 15626338
$ws.Range("M31").Value = -15626043
$ws.Range("H34").Value = 13700560
$ws.Range("I34").Value = 15626338
$ws.Range("K34").Value = 15626338
$ws.Range("M34").Value = -15626136
$ws.Range("H99").Value = 9272.272000000001
$ws.Range("I99").Value = 6665.8335
$ws.Range("J99").Value = 12400
$ws.Range("K99").Value = 6665.8335
$ws.Range("L99").Value = 12400
$ws.Range("M99").Value = -5167.8335
$ws.Range("N99").Value = -15396
$ws.Range("H113").Value = 919.6111
$ws.Range("I113").Value = 686.5
$ws.Range("J113").Value = 1385.8334
$ws.Range("K113").Value = 686.5
$ws.Range("L113").Value = 1385.8334
$ws.Range("M113").Value = 1483.5
$ws.Range("N113").Value = -5725.8334
$ws.Range("H115").Value = 49689.23
$ws.Range("J115").Value = 49689.23
$ws.Range("L115").Value = 49689.23
$ws.Range("N115").Value = -52039.23
$ws.Range("H126").Value = 9272.272000000001
$ws.Range("I126").Value = 6665.8335
$ws.Range("J126").Value = 12400
$ws.Range("K126").Value = 19997.5005
$ws.Range("L126").Value = 37200
$ws.Range("M126").Value = -17527.5005
$ws.Range("N126").Value = -42140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 749.0417
$ws.Range("J5").Value = 1015
$ws.Range("L5").Value = 3045
$ws.Range("N5").Value = -3269
$ws.Range("H38").Value = 172.6875
$ws.Range("I38").Value = 102.55556
$ws.Range("J38").Value = 262.85715
$ws.Range("K38").Value = 307.66668
$ws.Range("L38").Value = 788.5714499999999
$ws.Range("M38").Value = 39.33332000000001
$ws.Range("N38").Value = -1482.57145
$ws.Range("H86").Value = 324.66666
$ws.Range("J86").Value = 475
$ws.Range("L86").Value = 1425
$ws.Range("N86").Value = -3797
$ws.Range("H89").Value = 324.66666
$ws.Range("J89").Value = 475
$ws.Range("L89").Value = 4275
$ws.Range("N89").Value = -16131
$ws.Range("H107").Value = 532.3333
$ws.Range("J107").Value = 768.2
$ws.Range("L107").Value = 2304.6
$ws.Range("N107").Value = -6144.6
$ws.Range("H135").Value = 749.0417
$ws.Range("J135").Value = 1015
$ws.Range("L135").Value = 9135
$ws.Range("N135").Value = -14205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 49710
$ws.Range("J106").Value = 49710
$ws.Range("L106").Value = 49710
$ws.Range("N106").Value = -52234
$ws.Range("H114").Value = 43722
$ws.Range("J114").Value = 43722
$ws.Range("L114").Value = 43722
$ws.Range("N114").Value = -52400
$ws.Range("H118").Value = 25000
$ws.Range("I118").Value = 15000
$ws.Range("K118").Value = 15000
$ws.Range("M118").Value = -13343
$ws.Range("H132").Value = 3196.5957
$ws.Range("I132").Value = 2756.2703
$ws.Range("K132").Value = 8268.8109
$ws.Range("M132").Value = -5738.8109
$ws.Range("H139").Value = 80828.5
$ws.Range("J139").Value = 80828.5
$ws.Range("L139").Value = 80828.5
$ws.Range("N139").Value = -91108.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 726.44446
$ws.Range("I22").Value = 686.9167
$ws.Range("J22").Value = 758.06665
$ws.Range("K22").Value = 686.9167
$ws.Range("L22").Value = 758.06665
$ws.Range("M22").Value = -391.9167
$ws.Range("N22").Value = -1348.06665
$ws.Range("H27").Value = 726.44446
$ws.Range("I27").Value = 686.9167
$ws.Range("J27").Value = 758.06665
$ws.Range("K27").Value = 686.9167
$ws.Range("L27").Value = 758.06665
$ws.Range("M27").Value = -579.9167
$ws.Range("N27").Value = -972.06665
$ws.Range("H46").Value = 3438.8
$ws.Range("I46").Value = 2554.889
$ws.Range("J46").Value = 3817.6191
$ws.Range("K46").Value = 2554.889
$ws.Range("L46").Value = 3817.6191
$ws.Range("M46").Value = -2366.889
$ws.Range("N46").Value = -4193.6191
$ws.Range("H93").Value = 1151.7646
$ws.Range("I93").Value = 1042.5769
$ws.Range("K93").Value = 1042.5769
$ws.Range("M93").Value = 205.4231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 47499
$ws.Range("I61").Value = 50000
$ws.Range("J61").Value = 44998
$ws.Range("K61").Value = 50000
$ws.Range("L61").Value = 44998
$ws.Range("M61").Value = -49708
$ws.Range("N61").Value = -45582
$ws.Range("H113").Value = 781
$ws.Range("I113").Value = 708.8889
$ws.Range("K113").Value = 2126.6667
$ws.Range("M113").Value = 43.33329999999978
$ws.Range("H132").Value = 1792.3823
$ws.Range("I132").Value = 1239.931
$ws.Range("K132").Value = 3719.793
$ws.Range("M132").Value = -1189.793
$ws.Range("H136").Value = 4509.357
$ws.Range("I136").Value = 3204.15
$ws.Range("J136").Value = 7772.375
$ws.Range("K136").Value = 9612.450000000001
$ws.Range("L136").Value = 23317.125
$ws.Range("M136").Value = -7062.450000000001
$ws.Range("N136").Value = -28417.125

